$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''247.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''0.80%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''29.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''8.51%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.179'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''1.28%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.05734'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''0.79%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''6.575'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''0.74%'
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.095'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''2.85%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.8578'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''4.67%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = '''0.8696'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''1.64%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1365'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''2.25%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07064'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''1.76%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.02912'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''1.73%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09386'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''0.01%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001525'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''0.29%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = '''0.04124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''1.21%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '''0.0006007'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-0.40%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.005961'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''-4.10%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D18").Value = '''0.007489'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''5,067.81%'
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '''3.490'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''-0.61%'
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''2.190'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''-1.80%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.3186'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''0.67%'
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''5.68%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.1286'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''0.93%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''3.464'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-2.57%'
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''0.38%'
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''12.02%'
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''0.51%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.0001210'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''22.17%'
$ws.Range("E28").Style = "Normal"
$ws.Range("E40").Value = '''0.79%'
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1073'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''1.48%'
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002449'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''0.37%'
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003482'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''-42.07%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.008469'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-12.86%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.00005252'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''5.01%'
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''-0.06%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''-35.98%'
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''-9.30%'
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''-0.06%'
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''-0.06%'
$ws.Range("E50").Style = "Normal"
